# Update "想去人数" (F column) figures across sheets, per commit
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 5491
$ws.Cells.Item(5, 6).Value = 5491
$ws.Cells.Item(6, 6).Value = 215
$ws.Cells.Item(9, 6).Value = 1238
$ws.Cells.Item(12, 6).Value = 807
$ws.Cells.Item(13, 6).Value = 22
$ws.Cells.Item(14, 6).Value = 6547
$ws.Cells.Item(15, 6).Value = 39
$ws.Cells.Item(16, 6).Value = 87
$ws.Cells.Item(17, 6).Value = 123
$ws.Cells.Item(18, 6).Value = 4763
$ws.Cells.Item(19, 6).Value = 112
$ws.Cells.Item(21, 6).Value = 4195
$ws.Cells.Item(23, 6).Value = 4125
$ws.Cells.Item(25, 6).Value = 212
$ws.Cells.Item(26, 6).Value = 282
$ws.Cells.Item(27, 6).Value = 265
$ws.Cells.Item(29, 6).Value = 121
$ws.Cells.Item(31, 6).Value = 58
$ws.Cells.Item(32, 6).Value = 144
$ws.Cells.Item(33, 6).Value = 58
$ws.Cells.Item(34, 6).Value = 7467
$ws.Cells.Item(36, 6).Value = 1254
$ws.Cells.Item(37, 6).Value = 617
$ws.Cells.Item(38, 6).Value = 114
$ws.Cells.Item(39, 6).Value = 981
$ws.Cells.Item(41, 6).Value = 1486
$ws.Cells.Item(42, 6).Value = 195
$ws.Cells.Item(43, 6).Value = 828
$ws.Cells.Item(45, 6).Value = 3639
$ws.Cells.Item(46, 6).Value = 333
$ws.Cells.Item(47, 6).Value = 14
$ws.Cells.Item(49, 6).Value = 815
$ws.Cells.Item(50, 6).Value = 1027

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6, 6).Value = 69
$ws.Cells.Item(7, 6).Value = 8
$ws.Cells.Item(10, 6).Value = 8
$ws.Cells.Item(12, 6).Value = 23
$ws.Cells.Item(18, 6).Value = 66
$ws.Cells.Item(21, 6).Value = 849

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 229

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 229
$ws.Cells.Item(6, 6).Value = 5491
$ws.Cells.Item(7, 6).Value = 5491
$ws.Cells.Item(8, 6).Value = 215
$ws.Cells.Item(12, 6).Value = 1238
$ws.Cells.Item(15, 6).Value = 8
$ws.Cells.Item(16, 6).Value = 807
$ws.Cells.Item(17, 6).Value = 6547
$ws.Cells.Item(18, 6).Value = 39
$ws.Cells.Item(19, 6).Value = 87
$ws.Cells.Item(20, 6).Value = 123
$ws.Cells.Item(21, 6).Value = 4763
$ws.Cells.Item(22, 6).Value = 112
$ws.Cells.Item(24, 6).Value = 4195
$ws.Cells.Item(25, 6).Value = 4125
$ws.Cells.Item(27, 6).Value = 212
$ws.Cells.Item(28, 6).Value = 282
$ws.Cells.Item(31, 6).Value = 121
$ws.Cells.Item(32, 6).Value = 58
$ws.Cells.Item(34, 6).Value = 7467
$ws.Cells.Item(36, 6).Value = 1254
$ws.Cells.Item(37, 6).Value = 617
$ws.Cells.Item(38, 6).Value = 114
$ws.Cells.Item(39, 6).Value = 981
$ws.Cells.Item(41, 6).Value = 1486
$ws.Cells.Item(42, 6).Value = 195
$ws.Cells.Item(43, 6).Value = 828
$ws.Cells.Item(45, 6).Value = 3639
$ws.Cells.Item(46, 6).Value = 333
$ws.Cells.Item(48, 6).Value = 815
$ws.Cells.Item(49, 6).Value = 1027

Write-Output "Updated F-column counts across 展览, 演出, 本地生活, 全部类型 sheets."
